$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$dateFmt = "yyyy\-mm\-dd;@"

# Row 37: add Midterm 1 exam + date (matches C37 date)
$ws.Range("D37").Value = "Midterm 1"
$ws.Range("E37").Value = $ws.Range("C37").Value2
$ws.Range("E37").NumberFormat = $dateFmt

# Row 38: remove Midterm 1 exam + date (moved to row 37)
$ws.Range("D38").Clear()
$ws.Range("E38").Clear()

# Row 42: remove Midterm 2 exam + date (moved to row 44)
$ws.Range("D42").Clear()
$ws.Range("E42").Clear()

# Row 44: add Midterm 2 exam + date (matches C44 date)
$ws.Range("D44").Value = "Midterm 2"
$ws.Range("E44").Value = $ws.Range("C44").Value2
$ws.Range("E44").NumberFormat = $dateFmt

# Row 45: remove Midterm 3 exam + date (moved to row 48)
$ws.Range("D45").Clear()
$ws.Range("E45").Clear()

# Row 48: add Midterm 3 exam + date (matches C48 date)
$ws.Range("D48").Value = "Midterm 3"
$ws.Range("E48").Value = $ws.Range("C48").Value2
$ws.Range("E48").NumberFormat = $dateFmt

# Update selection to E45
[void]$ws.Range("E45").Select()
